# API-Automated_testing V1.1: Added comments for all classes and methods
#
# This script reproduces, via the Excel COM object model, the edits that were
# made to tests/data/test_data.xlsx:
#   - RegisterData sheet: selection changed to A1:A4
#   - CheckUsernameData sheet: a new column was inserted before column A,
#     populated with the "username"/"a123456"/"a654321"/"a456789" values
#     (mirroring the data already present in the shifted-to-B column), the
#     new column was narrowed, and the sheet's selection moved to G14.

$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("RegisterData")
$wsCheck    = $wb.Worksheets.Item("CheckUsernameData")

# --- RegisterData: shrink the current selection down to column A ----------
$wsRegister.Range("A1:A4").Select()

# --- CheckUsernameData: insert a new column A ------------------------------
# The existing column A (username header, expectedResult values, styles and
# all) shifts one place to the right and becomes column B.
$wsCheck.Columns("A").Insert()

# Give the new column its own (narrower) width.
$wsCheck.Columns("A").ColumnWidth = 17.86

# Fill the new column A with the same data as column B (it mirrors the
# "username" column used elsewhere in the workbook).
$wsCheck.Range("A1").Value = "username"
$wsCheck.Range("A2").Value = "a123456"
$wsCheck.Range("A3").Value = "a654321"
$wsCheck.Range("A4").Value = "a456789"

# Also fill in B3/B4, which were previously blank placeholders, with the
# same values now shown in column A.
$wsCheck.Range("B3").Value = "a654321"
$wsCheck.Range("B4").Value = "a456789"

# Pick up the formatting (number format / alignment) used for the header and
# data rows elsewhere in the workbook so the new column matches the look of
# the existing one.
$wsRegister.Range("A1").Copy()
$wsCheck.Range("A1").PasteSpecial(-4122)
$wsRegister.Range("A2").Copy()
$wsCheck.Range("A2").PasteSpecial(-4122)
$wsRegister.Range("A3").Copy()
$wsCheck.Range("A3").PasteSpecial(-4122)
$wsRegister.Range("A4").Copy()
$wsCheck.Range("A4").PasteSpecial(-4122)

# Finally, move the active selection on this (active) sheet to G14.
$wsCheck.Range("G14").Select()
